$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Turn A1:U62 into an Excel Table ("Table1").
#    A table built directly on a pre-formatted header row (ours is bold/
#    shaded) makes Excel capture that formatting into a styles.xml <dxf> and
#    reference it via headerRowDxfId - which the target workbook does not
#    have. So the table is first created on a small, unformatted helper
#    range, then resized onto the real A1:U62 range (Resize keeps the table
#    free of any captured header dxf), and the helper cells are cleared
#    afterwards.
# ---------------------------------------------------------------------------
$ws.Range("AA1").Value = "h1"
$ws.Range("AB1").Value = "h2"
$ws.Range("AA2").Value = "a"
$ws.Range("AB2").Value = "b"

$helperRange = $ws.Range("AA1:AB2")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $helperRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.Resize($ws.Range("A1:U62"))
$lo.TableStyle = ""

# ---------------------------------------------------------------------------
# 2) Rename the header row shared strings: "_old" -> "_FV2210", "_new" -> "_FV2304"
#    Doing this after the table is (re)sized onto A1:U62 makes the table's
#    column names track the real header text.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Segmentname_FV2210"
$ws.Range("B1").Value = "Segmentgruppe_FV2210"
$ws.Range("C1").Value = "Segment_FV2210"
$ws.Range("D1").Value = "Datenelement_FV2210"
$ws.Range("E1").Value = "Segment ID_FV2210"
$ws.Range("F1").Value = "Code_FV2210"
$ws.Range("G1").Value = "Qualifier_FV2210"
$ws.Range("H1").Value = "Beschreibung_FV2210"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2210"
$ws.Range("J1").Value = "Bedingung_FV2210"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2304"
$ws.Range("M1").Value = "Segmentgruppe_FV2304"
$ws.Range("N1").Value = "Segment_FV2304"
$ws.Range("O1").Value = "Datenelement_FV2304"
$ws.Range("P1").Value = "Segment ID_FV2304"
$ws.Range("Q1").Value = "Code_FV2304"
$ws.Range("R1").Value = "Qualifier_FV2304"
$ws.Range("S1").Value = "Beschreibung_FV2304"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("U1").Value = "Bedingung_FV2304"

$ws.Range("AA1:AB2").Clear()

# ---------------------------------------------------------------------------
# 3) Freeze the header row (row 1) and keep the pane/selection on row 2.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
